# artwork_data_colores2.xlsx ("Artistas" sheet)
#
# The commit shuffles the order of several tied entries in the artist
# frequency table (rows whose column-B count is equal, e.g. the big block
# of count == 1 rows) — a re-sort of duplicate-count groups that isn't
# stable across runs. Column B (the counts) is untouched; only the column A
# artist-name labels move between rows. Rewrite column A for every row
# whose label changed so it matches the new row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Wols"
$ws.Range("A5").Value = "Jones, Allen"
$ws.Range("A12").Value = "Collins, Cecil"
$ws.Range("A13").Value = "Tyson, Ian"
$ws.Range("A19").Value = "Palmer, Roger"
$ws.Range("A21").Value = "Stephenson, Ian"
$ws.Range("A22").Value = "Nash, John"
$ws.Range("A27").Value = "Hodgkin, Howard"
$ws.Range("A28").Value = "Manessier, Alfred"
$ws.Range("A29").Value = "Salt, John"
$ws.Range("A30").Value = "Phillips, Esq Tom"
$ws.Range("A31").Value = "Cutts, Simon"
$ws.Range("A32").Value = "Art & Language (Michael Baldwin, born 1945; Mel Ramsden, born 1944)"
$ws.Range("A33").Value = "Pomodoro, Gio"
$ws.Range("A34").Value = "Baselitz, Georg"
$ws.Range("A35").Value = "Matisse, Henri"
$ws.Range("A36").Value = "Park, Alistair"
$ws.Range("A37").Value = "Le Parc, Julio"
$ws.Range("A38").Value = "Wentworth, Richard"
$ws.Range("A39").Value = "Blake, John"
$ws.Range("A40").Value = "Spencer, Sir Stanley"
$ws.Range("A41").Value = "Benjamin, Anthony"
$ws.Range("A42").Value = "Lindström, Bengt"
$ws.Range("A43").Value = "Schneider, Gerard"
$ws.Range("A44").Value = "Grayson, Roy"
$ws.Range("A45").Value = "Dine, Jim"
$ws.Range("A46").Value = "Fabro, Luciano"
$ws.Range("A47").Value = "London Gallery"
$ws.Range("A48").Value = "Ackroyd, Norman"
$ws.Range("A49").Value = "Baumeister, Willi"
$ws.Range("A51").Value = "Benrath, Frédéric"
$ws.Range("A52").Value = "Götz, Professor Karl-Otto"
$ws.Range("A53").Value = "King, Ronald"
$ws.Range("A54").Value = "Cuixart, Modestos"
$ws.Range("A55").Value = "Illes, Arpad"
$ws.Range("A56").Value = "Murphy, John"
$ws.Range("A57").Value = "Hamilton Finlay, Ian"
$ws.Range("A58").Value = "Bird, John"
$ws.Range("A59").Value = "Nevinson, Christopher Richard Wynne"
$ws.Range("A60").Value = "Braque, Georges"
$ws.Range("A61").Value = "Appel, Karel"
$ws.Range("A62").Value = "Maccari, Mino"
$ws.Range("A63").Value = "Loker, John"
$ws.Range("A64").Value = "Neiland, Brendan"
$ws.Range("A65").Value = "Disler, Martin"
$ws.Range("A66").Value = "Thomkins, André"
$ws.Range("A67").Value = "Rocamora, Jaume"
$ws.Range("A68").Value = "Rainer, Arnulf"
$ws.Range("A69").Value = "Soulages, Pierre"
$ws.Range("A70").Value = "Guston, Philip"
$ws.Range("A71").Value = "Abrahams, Ivor"
$ws.Range("A72").Value = "Downsbrough, Peter"
$ws.Range("A73").Value = "Kunkel, Don"
$ws.Range("A74").Value = "Hayter, Stanley William"
$ws.Range("A75").Value = "Wunderlich, Paul"
$ws.Range("A76").Value = "Herring, Ed"
$ws.Range("A77").Value = "Blake, Peter"
$ws.Range("A78").Value = "Frohner, Adolf"
$ws.Range("A79").Value = "Laabs, Hans"
$ws.Range("A80").Value = "Sutherland, Graham, OM"
$ws.Range("A81").Value = "Beuys, Joseph"
$ws.Range("A82").Value = "Irvin, Albert"
$ws.Range("A83").Value = "Tamayo, Rufino"
$ws.Range("A84").Value = "Lattanzi, Luciano"
$ws.Range("A85").Value = "Trevelyan, Julian"
